$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value into a cell while preventing Excel from
# auto-converting number-like or percent-like strings into numeric values,
# and then strip the temporary Text number-format so the cell keeps the
# workbook default (General) formatting.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '310.28'
Set-TextValue $ws.Range("E2") '-1.63%'
Set-TextValue $ws.Range("D3") '37.59'
Set-TextValue $ws.Range("E3") '-4.12%'
Set-TextValue $ws.Range("D4") '5.108'
Set-TextValue $ws.Range("E4") '-0.93%'
Set-TextValue $ws.Range("D5") '0.07755'
Set-TextValue $ws.Range("E5") '-4.99%'
Set-TextValue $ws.Range("D6") '4.348'
Set-TextValue $ws.Range("E6") '-1.32%'
Set-TextValue $ws.Range("D7") '8.203'
Set-TextValue $ws.Range("E7") '-1.77%'
Set-TextValue $ws.Range("D8") '1.877'
Set-TextValue $ws.Range("E8") '-4.87%'
Set-TextValue $ws.Range("E9") '-11.76%'
Set-TextValue $ws.Range("D10") '0.9196'
Set-TextValue $ws.Range("E10") '-1.73%'
Set-TextValue $ws.Range("D11") '0.1194'
Set-TextValue $ws.Range("E11") '-8.72%'
Set-TextValue $ws.Range("D12") '0.1912'
Set-TextValue $ws.Range("E12") '-3.66%'
Set-TextValue $ws.Range("D13") '0.08858'
Set-TextValue $ws.Range("E13") '-1.71%'
Set-TextValue $ws.Range("D14") '0.03395'
Set-TextValue $ws.Range("E14") '-2.66%'
Set-TextValue $ws.Range("D15") '0.09685'
Set-TextValue $ws.Range("E15") '-0.43%'
Set-TextValue $ws.Range("D16") '0.001374'
Set-TextValue $ws.Range("E16") '-2.53%'
Set-TextValue $ws.Range("D17") '0.005775'
Set-TextValue $ws.Range("E17") '-8.00%'
Set-TextValue $ws.Range("D18") '3.552'
Set-TextValue $ws.Range("E18") '-1.62%'
Set-TextValue $ws.Range("D19") '0.3407'
Set-TextValue $ws.Range("E19") '-1.79%'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range("D20") '5.035'
Set-TextValue $ws.Range("E20") '0.24%'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range("D21") '0.1262'
Set-TextValue $ws.Range("E21") '-3.63%'
Set-TextValue $ws.Range("D22") '0.2594'
Set-TextValue $ws.Range("E22") '4.08%'
Set-TextValue $ws.Range("D23") '0.02105'
Set-TextValue $ws.Range("E23") '5,591.59%'
Set-TextValue $ws.Range("D24") '0.04383'
Set-TextValue $ws.Range("E24") '-0.08%'
Set-TextValue $ws.Range("D25") '0.001213'
Set-TextValue $ws.Range("E25") '-2.74%'
Set-TextValue $ws.Range("D26") '0.004240'
Set-TextValue $ws.Range("E26") '-10.74%'
Set-TextValue $ws.Range("D27") '0.0001350'
Set-TextValue $ws.Range("E27") '-65.34%'
Set-TextValue $ws.Range("D39") '0.02082'
Set-TextValue $ws.Range("E39") '-6.96%'
Set-TextValue $ws.Range("D40") '0.04941'
Set-TextValue $ws.Range("E40") '-5.55%'
Set-TextValue $ws.Range("D41") '0.007675'
Set-TextValue $ws.Range("E41") '-0.93%'
Set-TextValue $ws.Range("D42") '0.009880'
Set-TextValue $ws.Range("E42") '-4.46%'
Set-TextValue $ws.Range("D43") '0.1340'
Set-TextValue $ws.Range("E43") '-4.13%'
Set-TextValue $ws.Range("D44") '0.002060'
Set-TextValue $ws.Range("E44") '-2.03%'
Set-TextValue $ws.Range("D45") '0.009629'
Set-TextValue $ws.Range("E45") '5.47%'
Set-TextValue $ws.Range("D46") '0.00006568'
Set-TextValue $ws.Range("E46") '-3.81%'
Set-TextValue $ws.Range("D47") '0.00000000750'
Set-TextValue $ws.Range("E47") '-0.14%'
Set-TextValue $ws.Range("D48") '0.003043'
Set-TextValue $ws.Range("E48") '1.09%'
Set-TextValue $ws.Range("D50") '0.00002100'
Set-TextValue $ws.Range("E50") '-0.14%'
Set-TextValue $ws.Range("D51") '0.0002000'
Set-TextValue $ws.Range("E51") '-0.14%'
